# "Add Leave Card 10/32023 3:18 PM"
# Applies corrections to the 2018 LEAVE CREDITS and 2017 LEAVE BALANCE
# sheets of the leave card workbook:
#   - Corrects the monthly PERIOD dates in "2018 LEAVE CREDITS" (2018-2022)
#     from month-start to month-end so they are consistent with the rest
#     of the table.
#   - Posts the EARNED leave (1.25/mo) for Aug-Dec 2023 and adds the
#     "2024" year-divider row, continuing the schedule into Jan-Mar 2024.
#   - Records two new Vacation Leave applications (10 and 22 days) on the
#     "2017 LEAVE BALANCE" sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws3 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# ---------------------------------------------------------------------
# 2017 LEAVE BALANCE: new VL applications (rows 22-24)
# ---------------------------------------------------------------------
$ws3.Range("A22").Value = 45139

$ws3.Range("A23").Value = 45170
$ws3.Range("B23").Value = "VL(10-0-0)"
$ws3.Range("D23").Value = 10
$ws3.Range("K23").Value = "9/18-22, 25-29/2023"

$ws3.Range("A24").Value = 45200
$ws3.Range("B24").Value = "VL(22-0-0)"
$ws3.Range("D24").Value = 22
$ws3.Range("K24").Value = "10/2-31/2023"

# ---------------------------------------------------------------------
# 2018 LEAVE CREDITS: correct PERIOD end-dates for 2018-2022 (month-start
# values were entered originally; replace with the month-end values used
# consistently from 2023 onward)
# ---------------------------------------------------------------------
$ws2.Range("A11").Value = 43131
$ws2.Range("A12").Value = 43159
$ws2.Range("A13").Value = 43190
$ws2.Range("A14").Value = 43220
$ws2.Range("A15").Value = 43251
$ws2.Range("A16").Value = 43281
$ws2.Range("A17").Value = 43312
$ws2.Range("A18").Value = 43343
$ws2.Range("A19").Value = 43373
$ws2.Range("A20").Value = 43404
$ws2.Range("A21").Value = 43434
$ws2.Range("A22").Value = 43465
$ws2.Range("A24").Value = 43496
$ws2.Range("A25").Value = 43524
$ws2.Range("A26").Value = 43555
$ws2.Range("A27").Value = 43585
$ws2.Range("A28").Value = 43616
$ws2.Range("A29").Value = 43646
$ws2.Range("A30").Value = 43677
$ws2.Range("A31").Value = 43708
$ws2.Range("A32").Value = 43738
$ws2.Range("A33").Value = 43769
$ws2.Range("A34").Value = 43799
$ws2.Range("A35").Value = 43830
$ws2.Range("A37").Value = 43861
$ws2.Range("A38").Value = 43890
$ws2.Range("A39").Value = 43921
$ws2.Range("A40").Value = 43951
$ws2.Range("A41").Value = 43982
$ws2.Range("A42").Value = 44012
$ws2.Range("A43").Value = 44043
$ws2.Range("A44").Value = 44074
$ws2.Range("A45").Value = 44104
$ws2.Range("A46").Value = 44135
$ws2.Range("A47").Value = 44165
$ws2.Range("A48").Value = 44196
$ws2.Range("A50").Value = 44227
$ws2.Range("A51").Value = 44255
$ws2.Range("A52").Value = 44286
$ws2.Range("A53").Value = 44316
$ws2.Range("A54").Value = 44347
$ws2.Range("A55").Value = 44377
$ws2.Range("A56").Value = 44408
$ws2.Range("A57").Value = 44439
$ws2.Range("A58").Value = 44469
$ws2.Range("A59").Value = 44500
$ws2.Range("A60").Value = 44530
$ws2.Range("A61").Value = 44561
$ws2.Range("A63").Value = 44592
$ws2.Range("A64").Value = 44620
$ws2.Range("A65").Value = 44651
$ws2.Range("A66").Value = 44681
$ws2.Range("A67").Value = 44712
$ws2.Range("A68").Value = 44742
$ws2.Range("A69").Value = 44773
$ws2.Range("A70").Value = 44804
$ws2.Range("A71").Value = 44834
$ws2.Range("A72").Value = 44865
$ws2.Range("A73").Value = 44895
$ws2.Range("A74").Value = 44926

# ---------------------------------------------------------------------
# 2018 LEAVE CREDITS: post EARNED leave credits for Aug-Dec 2023
# ---------------------------------------------------------------------
$ws2.Range("C83").Value = 1.25
$ws2.Range("C84").Value = 1.25
$ws2.Range("C85").Value = 1.25
$ws2.Range("C86").Value = 1.25

$ws2.Range("B87").Value = "FL(5-0-0)"
$ws2.Range("C87").Value = 1.25
$ws2.Range("D87").Value = 5

# Row 88 becomes the "2024" year-divider row (matches the style used by
# the 2019-2023 divider rows, e.g. A10, A23, A36 ...)
$ws2.Range("A10").Copy()
$ws2.Range("A88").PasteSpecial(-4122)
$ws2.Range("A88").Value = "2024"

# Jan-Mar 2024 shift down into rows 89-91 (after the new divider row)
$ws2.Range("A89").Value = 45322
$ws2.Range("A90").Value = 45351
$ws2.Range("A91").Value = 45382

# Rows 92-93 no longer hold dates (Apr/May 2024 rows are not yet due)
$ws2.Range("A92").ClearContents()
$ws2.Range("A93").ClearContents()

# ---------------------------------------------------------------------
# Leave the workbook on the "2018 LEAVE CREDITS" tab, matching the
# activeTab/tabSelected state saved with this edit.
# ---------------------------------------------------------------------
$ws3.Range("F22").Select()
$ws2.Activate()
$ws2.Range("I9").Select()
